# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "29.109.63"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "1.832.21"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6279"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07533"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2930"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07683"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").Value = "1.833.13"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.028"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6687"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009385"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.995"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").Value = "29.097.29"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "2.074.78"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "223.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.141"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.500"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05764"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.164"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.122"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.204"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.836"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7411"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.139"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.668"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.766"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.226.33"
$ws.Range("E39").Value = "  -1.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01780"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.492"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8937"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "1.975.99"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000121"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5091"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07551"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +12.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4068"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.009"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.37%  "
